$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.938.71'
$ws.Range("E2").Value = '  +0.73%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.766.66'
$ws.Range("E3").Value = '  -0.48%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '328.45'
$ws.Range("E5").Value = '  +0.34%  '

$ws.Range("E6").Value = '  -0.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4656'
$ws.Range("E7").Value = '  +1.43%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3518'
$ws.Range("E8").Value = '  -1.96%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '43.52'
$ws.Range("E9").Value = '  +3.99%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07381'
$ws.Range("E10").Value = '  -1.51%  '

$ws.Range("E11").Value = '  -1.63%  '

$ws.Range("E12").Value = '  +0.01%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.60'
$ws.Range("E13").Value = '  -1.14%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.004'
$ws.Range("E14").Value = '  -0.71%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.183'
$ws.Range("E15").Value = '  -0.61%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.766.03'
$ws.Range("E16").Value = '  -0.45%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.21'
$ws.Range("E17").Value = '  -1.70%  '

$ws.Range("E18").Value = '  -0.54%  '

$ws.Range("E20").Value = '  -0.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.91'
$ws.Range("E21").Value = '  -1.09%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.787'
$ws.Range("E22").Value = '  -0.41%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.982.20'
$ws.Range("E23").Value = '  +0.66%  '

$ws.Range("E24").Value = '  -1.80%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.156'
$ws.Range("E25").Value = '  +3.56%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.06'
$ws.Range("E26").Value = '  -0.38%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.969.89'
$ws.Range("E28").Value = '  -0.40%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.184'
$ws.Range("E29").Value = '  +0.63%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '122.99'
$ws.Range("E30").Value = '  -2.26%  '

$ws.Range("E31").Value = '  -1.83%  '

$ws.Range("E32").Value = '  +0.85%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.653'
$ws.Range("E33").Value = '  -0.49%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.548'
$ws.Range("E34").Value = '  +0.28%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.68'

$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06103'
$ws.Range("E36").Value = '  -1.74%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02260'
$ws.Range("E37").Value = '  -1.75%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2066'
$ws.Range("E38").Value = '  -1.11%  '

$ws.Range("E39").Value = '  -0.90%  '

$ws.Range("B40").Value = 'WEMIXTOKEN'
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.452'
$ws.Range("E40").Value = '  +4.36%  '

$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6153'
$ws.Range("E41").Value = '  -2.55%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.189'
$ws.Range("E42").Value = '  +0.06%  '

$ws.Range("E43").Value = '  -0.43%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.10'
$ws.Range("E44").Value = '  -0.93%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.736'
$ws.Range("E45").Value = '  -0.19%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5791'
$ws.Range("E46").Value = '  -1.83%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '123.70'
$ws.Range("E47").Value = '  +1.03%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.931'
$ws.Range("E48").Value = '  -1.03%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06820'
$ws.Range("E49").Value = '  -1.53%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.123'
$ws.Range("E50").Value = '  -1.41%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.01'
$ws.Range("E51").Value = '  -0.43%  '
